# Update the NGA HPV16 summary table with refreshed model-fit statistics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "-67.1(se=23.1)"
$ws.Range("C2").Value = "0.011(95% CI, 0.0096-0.013)"

$ws.Range("B3").Value = "-45.9(se=11.1)"
$ws.Range("C3").Value = "0.66(95% CI, 0.21-1)"
$ws.Range("F3").Value = "1.4(95% CI, 0.45-2)"

$ws.Range("B4").Value = "-56.4(se=16.9)"

$ws.Range("B5").Value = "-66.3(se=21.7)"

$ws.Range("B6").Value = "-27(se=5)"
